# Augmented_datasets upload and graph_visualization
#
# - Rename header "name" (col B) to "badge_name"
# - Add a new "description" column (col G) with a short human-readable
#   summary (issuer + skills + difficulty) for every badge row
# - Minor view/formatting touch-ups that came along with the re-save
#   (column widths for the newly-wide text columns, cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Cells.Item(1, 2).Value = "badge_name"
$ws.Cells.Item(1, 7).Value = "description"
# new header cell picks up the same bold/centered header style as the rest
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)

# --- New "description" column, row by row ----------------------------
$ws.Cells.Item(2, 7).Value = "Coursera: Intro to Python & data analysis. Beginner."
$ws.Cells.Item(3, 7).Value = "DataCamp: Learn Pandas & SQL basics. Intermediate."
$ws.Cells.Item(4, 7).Value = "Udacity: ML fundamentals with Python. Intermediate."
$ws.Cells.Item(5, 7).Value = "DeepLearning.AI: Advanced neural nets & TensorFlow. Advanced."
$ws.Cells.Item(6, 7).Value = "AI Policy: Overview of ethics & policy in AI. Advanced."
$ws.Cells.Item(7, 7).Value = "AWS: Intro to cloud computing. Beginner."
$ws.Cells.Item(8, 7).Value = "Udacity: Basics of DevOps with Kubernetes & CI/CD. Intermediate."
$ws.Cells.Item(9, 7).Value = "edX: Fundamentals of cybersecurity & networking. Beginner."
$ws.Cells.Item(10, 7).Value = "Offensive Security: Advanced ethical hacking techniques. Advanced."
$ws.Cells.Item(11, 7).Value = "Google: Core UI/UX design fundamentals. Beginner."

# --- Column widths for the newly-visible text -------------------------
$ws.Columns.Item(2).ColumnWidth = 31.125
$ws.Columns.Item(6).ColumnWidth = 16.75
$ws.Columns.Item(7).ColumnWidth = 60.875

# --- Leave the cursor parked where the author left it ------------------
$ws.Range("G15").Select()
